# Apply updated cryptocurrency price/volume data scraped on
# Thu Aug 31 19:14:57 UTC 2023 (GitHub Actions refresh of cryptos list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.422.42'
$ws.Range('E2').Value = '  -3.60%  '
$ws.Range('D3').Value = '1.665.75'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.52'
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5158'
$ws.Range('E6').Value = '  -3.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06462'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2573'
$ws.Range('E9').Value = '  -3.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.95'
$ws.Range('E10').Value = '  -4.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07673'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.336'
$ws.Range('E12').Value = '  -5.28%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.667.54'
$ws.Range('E13').Value = '  -2.93%  '
$ws.Range('D14').Value = '1.897.54'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5541'
$ws.Range('E15').Value = '  -3.37%  '
$ws.Range('D16').Value = '0.0₅8063'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.69'
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('D18').Value = '26.485.90'
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '210.31'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.422'
$ws.Range('E21').Value = '  -5.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.10'
$ws.Range('E22').Value = '  -3.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.886'
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.11'
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.734'
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1164'
$ws.Range('E27').Value = '  -4.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.008'
$ws.Range('E28').Value = '  -3.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.76'
$ws.Range('E29').Value = '  -3.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05235'
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.262'
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.376'
$ws.Range('E32').Value = '  -3.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.220'
$ws.Range('E33').Value = '  -6.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.574'
$ws.Range('E34').Value = '  -4.44%  '
$ws.Range('E35').Value = '  -4.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.373'
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9255'
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5734'
$ws.Range('E38').Value = '  -2.33%  '
$ws.Range('D39').Value = '1.155.21'
$ws.Range('E39').Value = '  +10.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01599'
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8497'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.645'
$ws.Range('E43').Value = '  -3.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.23'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').Value = '1.806.84'
$ws.Range('E45').Value = '  -2.54%  '
$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4495'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.96'
$ws.Range('E48').Value = '  -3.65%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.966'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05117'
$ws.Range('E51').Value = '  -2.54%  '
